$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the lead-in text before the first "Rocchio":
#    "...We set a threshold k, such that if fewer than k documents are
#    returned by the wordnet expansion, we will run the "
#    becomes
#    "...When wordnet is ineffective, we trigger "
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "We set a threshold k, such that if fewer than k documents are returned by the wordnet expansion, we will run the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "When wordnet is ineffective, we trigger ", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Replace the trailing text after the first "Rocchio":
#    " feedback to get more results."
#    becomes
#    " expansion. If WordNet returns fewer than two terms, the term
#    being searched is rare so we trigger Rocchio expansion."
#    A placeholder is used for the second "Rocchio" so it can be
#    isolated afterwards.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " feedback to get more results.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " expansion. If WordNet returns fewer than two terms, the term being searched is rare so we trigger RocchioPLACEHOLDER expansion.", 2) | Out-Null

$d.Content.Find.Execute(
    "RocchioPLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false,
    "Rocchio", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Move the _GoBack bookmark from the end of the Conclusion
#    paragraph to the end of this paragraph (right after the new
#    final sentence we just appended, before the paragraph mark).
# ------------------------------------------------------------------
$oldBm = $d.Content.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$paraRng = $d.Content
$paraRng.Find.Execute("Nevertheless, recall is very important", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraRng.Expand(4) | Out-Null
$pilcrow = $paraRng.End - 1

# Placing a zero-length bookmark directly at a paragraph-end position is
# unreliable, so grow a temporary one-character range there, wrap it
# with the bookmark, then delete the character, leaving a collapsed
# bookmark anchored at the correct spot.
$insertPoint = $d.Range($pilcrow, $pilcrow)
$insertPoint.InsertBefore("X")
$markerRange = $d.Range($pilcrow, $pilcrow + 1)
$markerRange.Bookmarks.Add("_GoBack") | Out-Null
$newBm = $d.Content.Bookmarks.Item("_GoBack")
$newBm.Range.Text = ""
